$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape
$ws.Range("D2").Value = "62.236.56"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "3.435.15"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "411.66"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.47"
$ws.Range("E6").Value = "  -3.59%  "
$ws.Range("E7").Value = "  +5.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +12.70%  "
$ws.Range("E10").Value = "  +17.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.63"
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.74"
$ws.Range("E14").Value = "  +5.25%  "
$ws.Range("E15").Value = "  +56.40%  "
$ws.Range("D16").Value = "3.433.32"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.52"
$ws.Range("E17").Value = "  +13.52%  "
$ws.Range("E18").Value = "  +5.20%  "
$ws.Range("D19").Value = "62.284.98"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "403.77"
$ws.Range("E20").Value = "  +28.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "90.30"
$ws.Range("E21").Value = "  +8.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.20"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.48"
$ws.Range("E23").Value = "  +4.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.26"
$ws.Range("E24").Value = "  +3.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "33.66"
$ws.Range("E25").Value = "  +14.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.67"
$ws.Range("E26").Value = "  +4.23%  "
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.74"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.74"
$ws.Range("E29").Value = "  +10.07%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.172"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "44.07"
$ws.Range("E32").Value = "  +7.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.94"
$ws.Range("E33").Value = "  +5.26%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0504"
$ws.Range("E35").Value = "  +5.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "52.49"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.42"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.93"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("E40").Value = "  +6.89%  "
$ws.Range("E41").Value = "  +5.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.58"
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("E45").Value = "  +8.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.87"
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.13"
$ws.Range("E47").Value = "  +3.74%  "
$ws.Range("D48").Value = "2.126.08"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0373"
$ws.Range("E51").Value = "  +8.28%  "
